$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 'Permits Filed for 21–69 Shore Boulevard in Astoria, Queens'
$ws.Range("B2").Value = 'https://newyorkyimby.com/2025/10/permits-filed-for-21-69-shore-boulevard-in-astoria-queens.html'
$ws.Range("C2").Value = 'Permits have been filed for a four-story residential building at 21–69 Shore Boulevard in <a href="https://newyorkyimby.com/neighborhoods/astoria">Astoria</a>, Queens. Located between 21st Drive and Ditmars Boulevard, the lot is closest to the Astoria–Ditmars Boulevard subway station, served by the N and W trains. Aryeh Assouline of Impact Builders Corp. is listed as the owner behind the applications.'
$ws.Range("D2").Value = '2025-10-01T10:30:37+00:00'
$ws.Range("E2").Value = 'Wed, 01 Oct 2025 10:30:37 +0000'
$ws.Range("F2").Value = 'YIMBY'
$ws.Range("G2").Value = 'YIMBY - Astoria'
$ws.Range("H2").Value = ''
